# CryCompanywiseStockReport_1.xlsx — stock-quantity correction pass.
# For each affected line item, the quantity (col F) and its extended stock
# value (col G, = unit cost in col D * quantity) are corrected; the company
# "Sub Total:" rows (col B) and the final Sub Total / Grand Total rows (923/924)
# are recomputed to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 25: AL-AVLDAZLLER DEO N FRESH PERFUME BAG
$ws.Range("F25").Value = 41
$ws.Range("G25").Value = 1155.38

# Row 26: AL-AVLDAZLLER DEO N FRESH SANI CUBES DOUBLE
$ws.Range("F26").Value = 115
$ws.Range("G26").Value = 5301.5

# Row 34: AL-DAZLLER LIPCOLOUR ALL DAY 5g
$ws.Range("F34").Value = 39
$ws.Range("G34").Value = 3995.55

# Row 39: AL-DAZLLER NAIL POLISH REMOVER WIPES
$ws.Range("F39").Value = 39
$ws.Range("G39").Value = 1198.86

# Row 40: AL-DAZLLER WATER PROOF EYELINER 10ml
$ws.Range("F40").Value = 28
$ws.Range("G40").Value = 1290.8

# Row 46: Sub Total:
$ws.Range("B46").Value = 38543.23

# Row 48: ZOFF Almond 250GMS
$ws.Range("F48").Value = 110
$ws.Range("G48").Value = 21644.7

# Row 79: ZOFF Roasted Pistachio-250 GMS
$ws.Range("F79").Value = 0
$ws.Range("G79").Value = 0

# Row 85: Sub Total:
$ws.Range("B85").Value = 205808

# Row 123: BHA-Vicks inhaler
$ws.Range("F123").Value = 201
$ws.Range("G123").Value = 9024.9

# Row 134: Sub Total:
$ws.Range("B134").Value = 88580.60000000001

# Row 182: COL-Colgate Actice salt 100+200
$ws.Range("F182").Value = 4
$ws.Range("G182").Value = 573.76

# Row 186: COL-Colgate Kids 2+ Toothbrush
$ws.Range("F186").Value = 45
$ws.Range("G186").Value = 756.9

# Row 187: COL-Colgate Max Fresh Blue Toothpaste 150 gms
$ws.Range("F187").Value = 58
$ws.Range("G187").Value = 5387.62

# Row 198: Sub Total:
$ws.Range("B198").Value = 61200.69

# Row 224: DESAI-Mixed Fruit Jam 500g
$ws.Range("F224").Value = 81
$ws.Range("G224").Value = 6583.68

# Row 226: DESAI-Punjabi Masala Papad 160g
$ws.Range("F226").Value = 73
$ws.Range("G226").Value = 3336.1

# Row 228: Sub Total:
$ws.Range("B228").Value = 44941

# Row 256: GHP-Glamic Disinfactant surface cleaner 500ML
$ws.Range("F256").Value = 9
$ws.Range("G256").Value = 418.41

# Row 264: GHP-Nepthalene Balls 100 gms
$ws.Range("F264").Value = 52
$ws.Range("G264").Value = 1249.56

# Row 267: Sub Total:
$ws.Range("B267").Value = 35733.65

# Row 284: HAM-Exclutive Lunch Tiffin 3 Containers
$ws.Range("F284").Value = 27
$ws.Range("G284").Value = 7180.92

# Row 288: HAM-New Mug 1.5
$ws.Range("F288").Value = 172
$ws.Range("G288").Value = 4258.72

# Row 294: HAM-Thermosteel Duo Dlx 1000 Ml
$ws.Range("F294").Value = 26
$ws.Range("G294").Value = 21025.68

# Row 298: Sub Total:
$ws.Range("B298").Value = 137162.12

# Row 313: HIM-ALMOND & ROSE SOAP 125GX4NVALUE PAK
$ws.Range("F313").Value = 73
$ws.Range("G313").Value = 9894.42

# Row 319: HIM-AYURVEDA SANDAL GLOW SOAP 125G IND
$ws.Range("F319").Value = 293
$ws.Range("G319").Value = 9156.25

# Row 338: HIM-HIMALAYA AYURVEDA CLEAR SKIN SOAP 125G
$ws.Range("F338").Value = 118
$ws.Range("G338").Value = 4442.7

# Row 348: HIM-TOTAL CARE BABY PANTS DRAPERS-XL-9S
$ws.Range("F348").Value = 0
$ws.Range("G348").Value = 0

# Row 349: Sub Total:
$ws.Range("B349").Value = 159873.14

# Row 355: HUL-Cmft Fab Conditioner Ss Desire 850Ml
$ws.Range("F355").Value = 120
$ws.Range("G355").Value = 20734.8

# Row 389: Hul-pears pure and gentle 3x125 gm
$ws.Range("F389").Value = 82
$ws.Range("G389").Value = 10402.52

# Row 398: HUL-Ponds Pure Detox Fw 100G
$ws.Range("F398").Value = 9
$ws.Range("G398").Value = 1391.94

# Row 399: HUL-RAB 4x250g
$ws.Range("F399").Value = 358
$ws.Range("G399").Value = 21014.6

# Row 410: HUL-Surf Exl Mtc Liq Tl 1 Ltr Cp
$ws.Range("F410").Value = 57
$ws.Range("G410").Value = 8616.690000000001

# Row 422: HUL-Vim Liquid Yellow Bottle 750M
$ws.Range("F422").Value = 51
$ws.Range("G422").Value = 7326.66

# Row 423: Sub Total:
$ws.Range("B423").Value = 295979.17

# Row 436: HUL-Women Horlicks 400g Jar Chocolate
$ws.Range("F436").Value = 14
$ws.Range("G436").Value = 3130.4

# Row 437: Sub Total:
$ws.Range("B437").Value = 49759.77

# Row 460: JYOTHY-Exo Safai Small (Anti Bac)
$ws.Range("F460").Value = 27
$ws.Range("G460").Value = 179.55

# Row 464: JYOTHY-Mr. White Ultimate Whiteness Detergent Powder 3Kg (Bucket Free)
$ws.Range("F464").Value = 27
$ws.Range("G464").Value = 7060.77

# Row 479: JYT - Henko Matic Liquid Detergent Front Load 1L
$ws.Range("F479").Value = 37
$ws.Range("G479").Value = 5351.68

# Row 481: Sub Total:
$ws.Range("B481").Value = 65108.73

# Row 489: KAR-MYSORE SANDAL ROSE HANDWASH - 250 ML
$ws.Range("F489").Value = 3
$ws.Range("G489").Value = 165.99

# Row 492: KAR-MYSORE SANDAL TALC - 300 GM
$ws.Range("F492").Value = 19
$ws.Range("G492").Value = 1520.95

# Row 497: Sub Total:
$ws.Range("B497").Value = 61702.15

# Row 503: KCI-Huggies Wonderpants L64
$ws.Range("F503").Value = 5
$ws.Range("G503").Value = 3775

# Row 511: Sub Total:
$ws.Range("B511").Value = 41801.44

# Row 576: MAYA-Eva Rose 750W Mixer Grinder
$ws.Range("F576").Value = 3
$ws.Range("G576").Value = 8551.5

# Row 583: Sub Total:
$ws.Range("B583").Value = 101062.78

# Row 590: Ranga-Agb Baratanatyam 170g
$ws.Range("F590").Value = 9
$ws.Range("G590").Value = 456.03

# Row 592: RANGA-Cycle Brand Bansuri 75gms
$ws.Range("F592").Value = 194
$ws.Range("G592").Value = 6072.2

# Row 595: RANGA-Cycle Brand Morning Glory 105gms
$ws.Range("F595").Value = 121
$ws.Range("G595").Value = 6059.68

# Row 600: Ranga-Om Shanthi Navin Shinning Powder 200 Grams
$ws.Range("F600").Value = 91
$ws.Range("G600").Value = 2936.57

# Row 603: RANGA-Om Shanthi Saadana Puja Oil 1000ml Bottle
$ws.Range("F603").Value = 27
$ws.Range("G603").Value = 4301.64

# Row 605: RANGA-Stop - O Bathroom Freshner Power Bags Apple Cinnamon 10gms
$ws.Range("F605").Value = 107
$ws.Range("G605").Value = 3724.67

# Row 610: Sub Total:
$ws.Range("B610").Value = 87368.2

# Row 619: NES-Maggi 2 Minn Ndls Mas 280G
$ws.Range("F619").Value = 373
$ws.Range("G619").Value = 15994.24

# Row 621: NES-Maggi 2-Minn Special Mas 24X280G
$ws.Range("F621").Value = 313
$ws.Range("G621").Value = 19005.36

# Row 625: NES-Maggi Masala Noodles 420G
$ws.Range("F625").Value = 118
$ws.Range("G625").Value = 7589.76

# Row 638: Sub Total:
$ws.Range("B638").Value = 181106.88

# Row 669: PRI- B-06 KNIFE & PEELER 3 PC SET SS
$ws.Range("F669").Value = 63
$ws.Range("G669").Value = 4991.49

# Row 673: PRI-B 19 VIMAL Tea Strainer S.S
$ws.Range("F673").Value = 164
$ws.Range("G673").Value = 4277.12

# Row 684: PRI-B-58 VIMAL Fruit & Vegetable Push Chopper
$ws.Range("F684").Value = 47
$ws.Range("G684").Value = 8146.04

# Row 686: PRI-P-11 PRISTINE Paper Glass Disposable 250ML (Pack of 30pcs)
$ws.Range("F686").Value = 100
$ws.Range("G686").Value = 5448

# Row 688: Sub Total:
$ws.Range("B688").Value = 108165.04

# Row 779: SOU-Black Fard Dates 180g
$ws.Range("F779").Value = 16
$ws.Range("G779").Value = 1468

# Row 785: Sub Total:
$ws.Range("B785").Value = 22052.14

# Row 797: 24 Mantra SELECT ORGANIC KABULI CHANA.1 KG
$ws.Range("F797").Value = 4
$ws.Range("G797").Value = 674.64

# Row 800: 24 Mantra SELECT ORGANIC POHA-500 gms
$ws.Range("F800").Value = 24
$ws.Range("G800").Value = 1492.08

# Row 803: 24 Mantra SELECT ORGANIC RAJMA-1 KG
$ws.Range("F803").Value = 4
$ws.Range("G803").Value = 674.64

# Row 805: Sub Total:
$ws.Range("B805").Value = 8432.049999999999

# Row 813: SPI-Volini Spray 40gm
$ws.Range("F813").Value = 74
$ws.Range("G813").Value = 7439.22

# Row 814: SPI-Volini Spray 60 gm
$ws.Range("F814").Value = 118
$ws.Range("G814").Value = 16558.94

# Row 815: Sub Total:
$ws.Range("B815").Value = 47736.68

# Row 820: TCP-Kala Chana 1 kg
$ws.Range("F820").Value = 99
$ws.Range("G820").Value = 8074.44

# Row 823: TCP-Tata gold leaf 250g
$ws.Range("F823").Value = 193
$ws.Range("G823").Value = 21527.22

# Row 831: TCP-Tata Tea Chakra Gold 500 gm ( Tamilnadu / Andhra / Telangana)
$ws.Range("F831").Value = 58
$ws.Range("G831").Value = 15860.68

# Row 837: Sub Total:
$ws.Range("B837").Value = 250115.02

# Row 842: Orgfeed Arhar Dal 1 Kg
$ws.Range("F842").Value = 56
$ws.Range("G842").Value = 9177.84

# Row 859: Shankys Tip Top Gold Basmati Rice 1 Kg
$ws.Range("F859").Value = 374
$ws.Range("G859").Value = 38488.34

# Row 867: Sub Total:
$ws.Range("B867").Value = 274737

# Row 907: VVD Priyam Cold Pressed Groundnut Oil Pouch 1 Ltr
$ws.Range("F907").Value = 267
$ws.Range("G907").Value = 43550.37

# Row 911: Sub Total:
$ws.Range("B911").Value = 44515.7

# Row 923: Sub Total:
$ws.Range("B923").Value = 3932561.31

# Row 924: Note:Rates are Inclusive of Tax
$ws.Range("B924").Value = 3932561.31
